$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 911.125
$ws.Range("I62").Value = 927
$ws.Range("J62").Value = 800
$ws.Range("K62").Value = 927
$ws.Range("L62").Value = 800
$ws.Range("M62").Value = -303
$ws.Range("N62").Value = -2048

$ws.Range("H65").Value = 911.125
$ws.Range("I65").Value = 927
$ws.Range("J65").Value = 800
$ws.Range("K65").Value = 4635
$ws.Range("L65").Value = 4000
$ws.Range("M65").Value = -1515
$ws.Range("N65").Value = -10240

$ws.Range("H88").Value = 2488.875
$ws.Range("I88").Value = 1475.75
$ws.Range("J88").Value = 3502
$ws.Range("K88").Value = 1475.75
$ws.Range("L88").Value = 3502
$ws.Range("M88").Value = -1069.75
$ws.Range("N88").Value = -4314

$ws.Range("H91").Value = 2488.875
$ws.Range("I91").Value = 1475.75
$ws.Range("J91").Value = 3502
$ws.Range("K91").Value = 1475.75
$ws.Range("L91").Value = 3502
$ws.Range("M91").Value = -71.75
$ws.Range("N91").Value = -6310

$ws.Range("H134").Value = 52222.855
$ws.Range("J134").Value = 52222.855
$ws.Range("L134").Value = 52222.855
$ws.Range("N134").Value = -62362.855

$ws.Range("H137").Value = 2598.6182
$ws.Range("I137").Value = 2265.0667
$ws.Range("J137").Value = 4099.6
$ws.Range("K137").Value = 6795.2001
$ws.Range("L137").Value = 12298.8
$ws.Range("M137").Value = -4245.2001
$ws.Range("N137").Value = -17398.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4446.6934
$ws.Range("I32").Value = 3012.3704
$ws.Range("J32").Value = 14128.375
$ws.Range("K32").Value = 3012.3704
$ws.Range("L32").Value = 14128.375
$ws.Range("M32").Value = -2725.3704
$ws.Range("N32").Value = -14702.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4067.2727
$ws.Range("I105").Value = 4274
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 4274
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -2527
$ws.Range("N105").Value = -5494

$ws.Range("H134").Value = 4163.788
$ws.Range("I134").Value = 4444.517
$ws.Range("J134").Value = 2128.5
$ws.Range("K134").Value = 13333.551
$ws.Range("L134").Value = 6385.5
$ws.Range("M134").Value = -10798.551
$ws.Range("N134").Value = -11455.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32011.195
$ws.Range("I31").Value = 43101.42
$ws.Range("J31").Value = 3176.6
$ws.Range("K31").Value = 43101.42
$ws.Range("L31").Value = 3176.6
$ws.Range("M31").Value = -42806.42
$ws.Range("N31").Value = -3766.6

$ws.Range("H34").Value = 32011.195
$ws.Range("I34").Value = 43101.42
$ws.Range("J34").Value = 3176.6
$ws.Range("K34").Value = 43101.42
$ws.Range("L34").Value = 3176.6
$ws.Range("M34").Value = -42899.42
$ws.Range("N34").Value = -3580.6

$ws.Range("H58").Value = 1058.7667
$ws.Range("I58").Value = 1134.6
$ws.Range("J58").Value = 679.6
$ws.Range("K58").Value = 1134.6
$ws.Range("L58").Value = 679.6
$ws.Range("M58").Value = -931.5999999999999
$ws.Range("N58").Value = -1085.6

$ws.Range("H99").Value = 1709.2142
$ws.Range("I99").Value = 1633.4445
$ws.Range("J99").Value = 1845.6
$ws.Range("K99").Value = 1633.4445
$ws.Range("L99").Value = 1845.6
$ws.Range("M99").Value = -135.4445000000001
$ws.Range("N99").Value = -4841.6

$ws.Range("H126").Value = 1709.2142
$ws.Range("I126").Value = 1633.4445
$ws.Range("J126").Value = 1845.6
$ws.Range("K126").Value = 4900.333500000001
$ws.Range("L126").Value = 5536.799999999999
$ws.Range("M126").Value = -2430.333500000001
$ws.Range("N126").Value = -10476.8

$ws.Range("H134").Value = 9203.385
$ws.Range("I134").Value = 11102.4
$ws.Range("K134").Value = 33307.2
$ws.Range("M134").Value = -30772.2

$ws.Range("H136").Value = 1058.7667
$ws.Range("I136").Value = 1134.6
$ws.Range("J136").Value = 679.6
$ws.Range("K136").Value = 3403.8
$ws.Range("L136").Value = 2038.8
$ws.Range("M136").Value = -853.7999999999997
$ws.Range("N136").Value = -7138.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 862.97
$ws.Range("I131").Value = 290
$ws.Range("J131").Value = 886.84375
$ws.Range("K131").Value = 870
$ws.Range("L131").Value = 2660.53125
$ws.Range("M131").Value = 4170
$ws.Range("N131").Value = -12740.53125

$ws.Range("H138").Value = 3142.6
$ws.Range("I138").Value = 990
$ws.Range("J138").Value = 4577.6665
$ws.Range("K138").Value = 2970
$ws.Range("L138").Value = 13732.9995
$ws.Range("M138").Value = 2170
$ws.Range("N138").Value = -24012.9995

$ws.Range("H139").Value = 2692.4
$ws.Range("I139").Value = 931.125
$ws.Range("J139").Value = 3866.5833
$ws.Range("K139").Value = 2793.375
$ws.Range("L139").Value = 11599.7499
$ws.Range("M139").Value = 2346.625
$ws.Range("N139").Value = -21879.7499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 39750
$ws.Range("I35").Value = 29000
$ws.Range("J35").Value = 43333.332
$ws.Range("K35").Value = 29000
$ws.Range("L35").Value = 43333.332
$ws.Range("M35").Value = -28702
$ws.Range("N35").Value = -43929.332

$ws.Range("H43").Value = 9400.846
$ws.Range("I43").Value = 2926.375
$ws.Range("J43").Value = 19760
$ws.Range("K43").Value = 2926.375
$ws.Range("L43").Value = 19760
$ws.Range("M43").Value = -2775.375
$ws.Range("N43").Value = -20062

$ws.Range("H70").Value = 5017.423
$ws.Range("I70").Value = 5017.8
$ws.Range("J70").Value = 5016.909
$ws.Range("K70").Value = 5017.8
$ws.Range("L70").Value = 5016.909
$ws.Range("M70").Value = -4747.8
$ws.Range("N70").Value = -5556.909

$ws.Range("H73").Value = 5017.423
$ws.Range("I73").Value = 5017.8
$ws.Range("J73").Value = 5016.909
$ws.Range("K73").Value = 5017.8
$ws.Range("L73").Value = 5016.909
$ws.Range("M73").Value = -4081.8
$ws.Range("N73").Value = -6888.909

$ws.Range("H80").Value = 3004.7778
$ws.Range("I80").Value = 2711.25
$ws.Range("J80").Value = 3239.6
$ws.Range("K80").Value = 2711.25
$ws.Range("L80").Value = 3239.6
$ws.Range("M80").Value = -1713.25
$ws.Range("N80").Value = -5235.6

$ws.Range("H83").Value = 3004.7778
$ws.Range("I83").Value = 2711.25
$ws.Range("J83").Value = 3239.6
$ws.Range("K83").Value = 13556.25
$ws.Range("L83").Value = 16198
$ws.Range("M83").Value = -8564.25
$ws.Range("N83").Value = -26182

$ws.Range("H113").Value = 1630.1428
$ws.Range("I113").Value = 1630.1428
$ws.Range("K113").Value = 1630.1428
$ws.Range("M113").Value = 539.8571999999999

$ws.Range("H122").Value = 1450
$ws.Range("I122").Value = 1506
$ws.Range("J122").Value = 1268
$ws.Range("K122").Value = 4518
$ws.Range("L122").Value = 3804
$ws.Range("M122").Value = -2068
$ws.Range("N122").Value = -8704

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 21444.8
$ws.Range("I132").Value = 12716.444
$ws.Range("J132").Value = 100000
$ws.Range("K132").Value = 38149.33199999999
$ws.Range("L132").Value = 300000
$ws.Range("M132").Value = -35619.33199999999
$ws.Range("N132").Value = -305060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6875.75
$ws.Range("I122").Value = 7429.4287
$ws.Range("K122").Value = 22288.2861
$ws.Range("M122").Value = -19838.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 19714.5
$ws.Range("J123").Value = 19714.5
$ws.Range("L123").Value = 19714.5
$ws.Range("N123").Value = -29514.5
